$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$data = @{
    2  = @(0.8729054620558272, 1.815219781561128, 10.48527348284812, 3.238097200957395, 3.156020188197708, 42)
    3  = @(0.5705085376199325, 2.549581372391033, 17.96663688838797, 4.238706983077265, 4.252315348056046, 41)
    4  = @(0.9809733857494989, 2.755143531496009, 17.33767364708628, 4.16385322112659, 4.098200338427169, 40)
    5  = @(0.7796923468595006, 2.538240234672108, 17.06363191478983, 4.130814921391399, 4.10959326519414, 39)
    6  = @(0.9784690441693839, 3.157245746473258, 23.19647943196707, 4.816272358574323, 4.779135268779009, 38)
    7  = @(0.9337797391493523, 2.855712669763669, 18.20336195286742, 4.26653981029914, 4.220526698338592, 37)
    8  = @(0.9577349491508317, 3.053928057267465, 21.15389011660629, 4.599335834292414, 4.562326536244196, 36)
    9  = @(0.7951063176841026, 2.533398893255768, 16.69223698970055, 4.085613416575356, 4.066005058678526, 35)
    10 = @(0.7595042420159662, 2.834380573090335, 17.16476161472717, 4.143037727890873, 4.134075351224831, 34)
    11 = @(0.9752670721986967, 2.879719984736652, 21.52816472569717, 4.639845334243069, 4.606522625915246, 33)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}

$wb.Save()
